$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The original sheet has columns:
#   A categories | B title | C description | D input_output_description |
#   E std_input  | F std_output | G sample_code | H created_at | I modified_at |
#   J template   | K language (numeric code)
#
# and 9 data rows (rows 2-10).
#
# The target sheet:
#   - drops the two rows whose titles are "Bmi..." (old row 6) and
#     "輸出 n 行數字三角形" (old row 8)
#   - inserts a new column A "language" (text values: Python/C/Java) - the
#     old "categories" numeric column A is overwritten with the new language
#     values
#   - inserts a new column B "categories" (text values, e.g. 靜態方法/數學運算)
#     shifting former columns B..K (title..template) to C..K's neighbours
#   - removes the old numeric "language" column (formerly K)
# ---------------------------------------------------------------------------

# 1. Remove the two rows that are dropped entirely. Delete the
#    higher-numbered row first so the other row index stays valid.
$ws.Rows("8:8").Delete()
$ws.Rows("6:6").Delete()

# 2. Insert a new blank column at B. This shifts former columns B..K
#    (title .. language) one column to the right (C..L).
$ws.Columns("B:B").Insert()

# 3. The old numeric "language" column (originally K) is now at column L.
#    Drop it completely - it is not part of the target layout.
$ws.Columns("L:L").Delete()

# 4. Write the new header row.
$ws.Range("A1").Value = "language"
$ws.Range("B1").Value = "categories"

# 5. Overwrite column A with the textual language values, and fill in the
#    newly inserted column B with the category values, row by row.
$ws.Range("A2").Value = "Python"
$ws.Range("B2").Value = "靜態方法"

$ws.Range("A3").Value = "C"
$ws.Range("B3").Value = "靜態方法"

$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = "靜態方法"

$ws.Range("A5").Value = "Python"
$ws.Range("B5").Value = "靜態方法"

$ws.Range("A6").Value = "Java"
$ws.Range("B6").Value = "靜態方法"

$ws.Range("A7").Value = "Java"
$ws.Range("B7").Value = "數學運算"

$ws.Range("A8").Value = "Java"
$ws.Range("B8").Value = "數學運算"
